$wb = $excel.ActiveWorkbook
$wsFree = $wb.Worksheets.Item("FreeGames")
$ws = $wb.Worksheets.Item("Sheet1")

$ws.Cells.Item(1,1).Value = "Game Title"
$ws.Cells.Item(1,2).Value = "Add"
$ws.Cells.Item(1,3).Value = "Install"
$ws.Cells.Item(1,4).Value = "Status"

$ws.Cells.Item(2,1).Value = "King of the Board"
$ws.Cells.Item(2,2).Value = "No"
$ws.Cells.Item(2,3).Value = "No"
$ws.Cells.Item(3,1).Value = "Witches x Warlocks"
$ws.Cells.Item(3,2).Value = "No"
$ws.Cells.Item(3,3).Value = "No"
$ws.Cells.Item(4,1).Value = "Sector's Edge"
$ws.Cells.Item(4,2).Value = "No"
$ws.Cells.Item(4,3).Value = "No"
$ws.Cells.Item(5,1).Value = "仙剑奇侠传九野"
$ws.Cells.Item(5,2).Value = "No"
$ws.Cells.Item(5,3).Value = "No"
$ws.Cells.Item(6,1).Value = "Siren Head: Awakening"
$ws.Cells.Item(6,2).Value = "No"
$ws.Cells.Item(6,3).Value = "No"
$ws.Cells.Item(7,1).Value = "Happy's Humble Burger Farm Alpha"
$ws.Cells.Item(7,2).Value = "No"
$ws.Cells.Item(7,3).Value = "No"
$ws.Cells.Item(8,1).Value = "The Last Spell: Prologue"
$ws.Cells.Item(8,2).Value = "No"
$ws.Cells.Item(8,3).Value = "No"
$ws.Cells.Item(9,1).Value = "LA Monsters"
$ws.Cells.Item(9,2).Value = "No"
$ws.Cells.Item(9,3).Value = "No"
$ws.Cells.Item(10,1).Value = "GOD OF FLAME"
$ws.Cells.Item(10,2).Value = "No"
$ws.Cells.Item(10,3).Value = "No"
$ws.Cells.Item(11,1).Value = "Farmer And Tree"
$ws.Cells.Item(11,2).Value = "No"
$ws.Cells.Item(11,3).Value = "No"
$ws.Cells.Item(12,1).Value = "World of Soccer RELOADED"
$ws.Cells.Item(12,2).Value = "No"
$ws.Cells.Item(12,3).Value = "No"
$ws.Cells.Item(13,1).Value = "Warlock Bentspine - Toilet Edition"
$ws.Cells.Item(13,2).Value = "No"
$ws.Cells.Item(13,3).Value = "No"
$ws.Cells.Item(14,1).Value = "Greed Knights"
$ws.Cells.Item(14,2).Value = "No"
$ws.Cells.Item(14,3).Value = "No"
$ws.Cells.Item(15,1).Value = "Card Blitz: WWII"
$ws.Cells.Item(15,2).Value = "No"
$ws.Cells.Item(15,3).Value = "No"
$ws.Cells.Item(16,1).Value = "Ratten Reich - Dance of Kings"
$ws.Cells.Item(16,2).Value = "No"
$ws.Cells.Item(16,3).Value = "No"
$ws.Cells.Item(17,1).Value = "Oasis VR"
$ws.Cells.Item(17,2).Value = "No"
$ws.Cells.Item(17,3).Value = "No"
$ws.Cells.Item(18,1).Value = "Shotgun Witch"
$ws.Cells.Item(18,2).Value = "No"
$ws.Cells.Item(18,3).Value = "No"
$ws.Cells.Item(19,1).Value = "MannaRites"
$ws.Cells.Item(19,2).Value = "No"
$ws.Cells.Item(19,3).Value = "No"
$ws.Cells.Item(20,1).Value = "Tree Trunk Brook"
$ws.Cells.Item(20,2).Value = "No"
$ws.Cells.Item(20,3).Value = "No"
$ws.Cells.Item(21,1).Value = "Chaos Combat Chess"
$ws.Cells.Item(21,2).Value = "No"
$ws.Cells.Item(21,3).Value = "No"
$ws.Cells.Item(22,1).Value = "Aimi"
$ws.Cells.Item(22,2).Value = "No"
$ws.Cells.Item(22,3).Value = "No"
$ws.Cells.Item(23,1).Value = "Zero IDLE"
$ws.Cells.Item(23,2).Value = "No"
$ws.Cells.Item(23,3).Value = "No"
$ws.Cells.Item(24,1).Value = "Scribble It!"
$ws.Cells.Item(24,2).Value = "No"
$ws.Cells.Item(24,3).Value = "No"
$ws.Cells.Item(25,1).Value = "Blood of Steel"
$ws.Cells.Item(25,2).Value = "No"
$ws.Cells.Item(25,3).Value = "No"
$ws.Cells.Item(26,1).Value = "Circle of Sumo: Online Rumble!"
$ws.Cells.Item(26,2).Value = "No"
$ws.Cells.Item(26,3).Value = "No"
$ws.Cells.Item(27,1).Value = "WKSP Rumble"
$ws.Cells.Item(27,2).Value = "No"
$ws.Cells.Item(27,3).Value = "No"
$ws.Cells.Item(28,1).Value = "VR Only Binaural Odyssey"
$ws.Cells.Item(28,2).Value = "No"
$ws.Cells.Item(28,3).Value = "No"
$ws.Cells.Item(29,1).Value = "Vecter"
$ws.Cells.Item(29,2).Value = "No"
$ws.Cells.Item(29,3).Value = "No"
$ws.Cells.Item(30,1).Value = "RuneScape ®"
$ws.Cells.Item(30,2).Value = "No"
$ws.Cells.Item(30,3).Value = "No"
$ws.Cells.Item(31,1).Value = "Eternal Return: Black Survival"
$ws.Cells.Item(31,2).Value = "No"
$ws.Cells.Item(31,3).Value = "No"

$wsFree.Range("A1").Copy()
$ws.Range("A1:D1").PasteSpecial(-4122)
$ws.Range("A2:C31").PasteSpecial(-4122)

Write-Output "done"